$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Qui" (D column) reassignments -----------------------------------
# D8: "Huang, Liu, Wu" -> "Wu"
$ws.Range("D8").Value = "Wu"
# D9: "Huang, Liu, Wu" -> "Liu"
$ws.Range("D9").Value = "Liu"
# D11: "Wu, Liu" -> "HUANG"
$ws.Range("D11").Value = "HUANG"
# D12: "Liu, Huang" -> "Liu"
$ws.Range("D12").Value = "Liu"
# D13: "Huang, Wu" -> "Wu"
$ws.Range("D13").Value = "Wu"

# --- Row 8: fill in start/end dates + test result ----------------------
# H4/I4 already carry the date-formatted style used elsewhere in the sheet;
# copy that formatting onto H8/I8 before writing the values so the cells
# pick up the same style (rather than Excel minting a brand new numFmt).
$ws.Range("H4").Copy($ws.Range("H8:I8"))
$ws.Range("H8").Value = "3/6/2018"
$ws.Range("I8").Value = "3/6/2018"
$ws.Range("J8").Value = "OK"

# --- Row 9: fill in the start date only --------------------------------
$ws.Range("H4").Copy($ws.Range("H9"))
$ws.Range("H9").Value = "3/6/2018"

# --- Update the active selection ---------------------------------------
$ws.Range("K14").Select() | Out-Null
